$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Archive the old draft: insert 5 fresh rows at the top of the data block
# (rows 2-6) for the new PsychoPy stimulus set, pushing the previous rows down.
$ws.Range("A2:B6").EntireRow.Insert()
$ws.Range("A2:B6").ClearFormats()

# Full target data for rows 2..31 (path, isCulprit) written explicitly so the
# row order exactly matches the new layout (including the Leemet-Bockler block,
# whose culprit row now leads the block like every other block).
$data = @(
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/24-Kobe-Bryant-(Lakers)-selected/$24-Kobe-Bryant.png', $true)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/24-Kobe-Bryant-(Lakers)-selected/9-Matt-Barnes.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/24-Kobe-Bryant-(Lakers)-selected/14-Troy-Murphy.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/24-Kobe-Bryant-(Lakers)-selected/15-Metta-World-Peace.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/24-Kobe-Bryant-(Lakers)-selected/16-Pau-Gasol.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/11-Leemet-Bockler-(BC-Kalev)-selected/$11-Leemet-Bockler.png', $true)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/11-Leemet-Bockler-(BC-Kalev)-selected/2-Kasper-Surorg.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/11-Leemet-Bockler-(BC-Kalev)-selected/33-Martin-Meiers.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/11-Leemet-Bockler-(BC-Kalev)-selected/9-Severi-Kaukiainen.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/11-Leemet-Bockler-(BC-Kalev)-selected/15-Anrijs-Miska.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/2-Cooper-Flagg-(Duke-Blue-Devils)-selected/$2-Cooper-Flagg.png', $true)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/2-Cooper-Flagg-(Duke-Blue-Devils)-selected/52-Stanley-Borden.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/2-Cooper-Flagg-(Duke-Blue-Devils)-selected/55-Spencer-Hubbard.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/2-Cooper-Flagg-(Duke-Blue-Devils)-selected/7-Kon-Knueppel.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/2-Cooper-Flagg-(Duke-Blue-Devils)-selected/8-Darren-Harris.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/23-Stephanie-Mavunga-(Poland)-selected/$23-Stephanie-Mavunga.png', $true)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/23-Stephanie-Mavunga-(Poland)-selected/12-Liliana-Banaszak.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/23-Stephanie-Mavunga-(Poland)-selected/13-Weronika-Gajda.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/23-Stephanie-Mavunga-(Poland)-selected/15-Klaudia-Gertchen.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/23-Stephanie-Mavunga-(Poland)-selected/4-Julia-Niemojewska.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Marcelinho-Huertas-(La-Laguna-Tenerife)-selected/$9-Marcelinho-Huertas.png', $true)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Marcelinho-Huertas-(La-Laguna-Tenerife)-selected/15-Joan-Sastre.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Marcelinho-Huertas-(La-Laguna-Tenerife)-selected/3-Jaime-Fernandez.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Marcelinho-Huertas-(La-Laguna-Tenerife)-selected/35-Fran-Guerra.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Marcelinho-Huertas-(La-Laguna-Tenerife)-selected/42-Aaron-Doornekamp.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Terezia-Palenikova-(Slovakia)-selected/$9-Terezia-Palenikova.png', $true)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Terezia-Palenikova-(Slovakia)-selected/15-Nikola-Kovacikova.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Terezia-Palenikova-(Slovakia)-selected/3-Ivana-Jakubcova.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Terezia-Palenikova-(Slovakia)-selected/4-Veronika-Remenarova.png', $false)
    ,@('C:\Users\diego\OneDrive\Documents\Coding\Neurotech\Neurocrime\Display/9-Terezia-Palenikova-(Slovakia)-selected/88-Natalia-Martiskova.png', $false)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

Write-Host "Done. UsedRange=" $ws.UsedRange.Address()
